$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.917.18'
$ws.Range("E2").Value = '''  +3.25%  '
$ws.Range("D3").Value = '''2.612.05'
$ws.Range("E3").Value = '''  +1.54%  '
$ws.Range("E4").Value = '''  +0.01%  '
$ws.Range("D5").Value = '''571.55'
$ws.Range("E5").Value = '''  +0.44%  '
$ws.Range("D6").Value = '''143.52'
$ws.Range("E6").Value = '''  +0.34%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = '''  -0.41%  '
$ws.Range("D8").Value = '''0.601'
$ws.Range("E8").Value = '''  +1.14%  '
$ws.Range("D9").Value = '''2.640.43'
$ws.Range("E9").Value = '''  +2.63%  '
$ws.Range("D10").Value = '''6.52'
$ws.Range("E10").Value = '''  -2.37%  '
$ws.Range("E11").Value = '''  +3.25%  '
$ws.Range("E12").Value = '''  -2.94%  '
$ws.Range("D13").Value = '''0.370'
$ws.Range("E13").Value = '''  +7.52%  '
$ws.Range("D14").Value = '''3.085.68'
$ws.Range("E14").Value = '''  +1.92%  '
$ws.Range("D15").Value = '''60.918.17'
$ws.Range("E15").Value = '''  +3.17%  '
$ws.Range("D16").Value = '''23.58'
$ws.Range("E16").Value = '''  +5.41%  '
$ws.Range("E17").Value = '''  +3.08%  '
$ws.Range("D18").Value = '''2.616.50'
$ws.Range("E18").Value = '''  +1.44%  '
$ws.Range("D19").Value = '''11.35'
$ws.Range("E19").Value = '''  +11.14%  '
$ws.Range("E20").Value = '''  +3.41%  '
$ws.Range("D21").Value = '''347.84'
$ws.Range("E21").Value = '''  +3.50%  '
$ws.Range("D22").Value = '''7.14'
$ws.Range("E22").Value = '''  +14.53%  '
$ws.Range("E23").Value = '''  +0.37%  '
$ws.Range("D24").Value = '''0.524'
$ws.Range("E24").Value = '''  +14.54%  '
$ws.Range("D25").Value = '''64.12'
$ws.Range("E25").Value = '''  -0.44%  '
$ws.Range("E26").Value = '''  -0.27%  '
$ws.Range("E27").Value = '''  -0.33%  '
$ws.Range("D28").Value = '''7.73'
$ws.Range("D29").Value = '''0.0₃0799'
$ws.Range("E29").Value = '''  +2.93%  '
$ws.Range("D30").Value = '''1.81'
$ws.Range("E30").Value = '''  +7.84%  '
$ws.Range("D31").Value = '''0.996'
$ws.Range("E31").Value = '''  -0.22%  '
$ws.Range("D32").Value = '''6.32'
$ws.Range("E32").Value = '''  +4.19%  '
$ws.Range("D33").Value = '''161.24'
$ws.Range("E33").Value = '''  +1.76%  '
$ws.Range("D34").Value = '''19.50'
$ws.Range("E34").Value = '''  +2.80%  '
$ws.Range("E35").Value = '''  +5.98%  '
$ws.Range("D36").Value = '''0.963'
$ws.Range("E36").Value = '''  +10.28%  '
$ws.Range("E37").Value = '''  +5.12%  '
$ws.Range("D38").Value = '''1.59'
$ws.Range("E38").Value = '''  +6.17%  '
$ws.Range("D39").Value = '''37.76'
$ws.Range("E39").Value = '''  +1.64%  '
$ws.Range("D40").Value = '''0.858'
$ws.Range("E40").Value = '''  -1.31%  '
$ws.Range("D41").Value = '''3.80'
$ws.Range("E41").Value = '''  +3.60%  '
$ws.Range("D42").Value = '''297.57'
$ws.Range("E42").Value = '''  +1.71%  '
$ws.Range("D43").Value = '''139.66'
$ws.Range("E43").Value = '''  +11.48%  '
$ws.Range("B44").Value = '''Stellar'
$ws.Range("C44").Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '''0.0988'
$ws.Range("E44").Value = '''  +1.23%  '
$ws.Range("B45").Value = '''FirstDigitalUSD'
$ws.Range("C45").Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''0.994'
$ws.Range("E45").Value = '''  -0.49%  '
$ws.Range("B46").Value = '''Mantle'
$ws.Range("C46").Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '''0.607'
$ws.Range("E46").Value = '''  +2.78%  '
$ws.Range("B47").Value = '''Hedera'
$ws.Range("C47").Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").Value = '''0.0553'
$ws.Range("E47").Value = '''  +3.27%  '
$ws.Range("D48").Value = '''0.0242'
$ws.Range("E48").Value = '''  +4.05%  '
$ws.Range("D49").Value = '''10.70'
$ws.Range("E49").Value = '''  +0.64%  '
$ws.Range("D50").Value = '''19.76'
$ws.Range("E50").Value = '''  +7.25%  '
$ws.Range("B51").Value = '''Maker'
$ws.Range("C51").Value = '''https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '''2.042.71'
$ws.Range("E51").Value = '''  +5.00%  '
